# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) worksheet (the 5th sheet) gains three new trailing
# columns, populated for every existing data row:
#   H = date              -> "2012-04-24"
#   I = legislator_name   -> "蔡其昌"
#   J = legislator_id     -> 1377

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Header row
$ws.Cells.Item(1, 8).Value2 = "date"
$ws.Cells.Item(1, 9).Value2 = "legislator_name"
$ws.Cells.Item(1, 10).Value2 = "legislator_id"

# Find the last populated data row (column A holds serial numbers 69,70,...)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Force column H to text so the date-like string "2012-04-24" is kept
# literally instead of being re-interpreted as a date serial number.
$ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8)).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value2 = "2012-04-24"
    $ws.Cells.Item($r, 9).Value2 = "蔡其昌"
    $ws.Cells.Item($r, 10).Value2 = 1377
}
